$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace a paragraph's *content* (all runs, leaving the <w:p> and its
# <w:pPr> attributes intact) with freshly built run XML. We do this by taking
# a Range that spans the paragraph text but stops *before* the paragraph mark
# (End - 1) and calling Range.InsertXML on it; the host replaces the content
# of that paragraph in place instead of inserting a sibling paragraph.
# ---------------------------------------------------------------------------
function Set-ParagraphRunsXml($para, [string]$runsXml) {
    $prng = $para.Range
    $crng = $d.Range($prng.Start, $prng.End - 1)
    $xml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $crng.InsertXML($xml)
}

# ---------------------------------------------------------------------------
# 1) DOI paragraph: collapse 3 runs (with proofErr markers) into a single run.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("Digital Object Identifier 10.1109/ACCESS.2017.Doi Number", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$doiPara = $rng.Paragraphs(1)
$doiRunsXml = '<w:r><w:t>Digital Object Identifier 10.1109/ACCESS.2017.Doi Number</w:t></w:r>'
Set-ParagraphRunsXml $doiPara $doiRunsXml

# ---------------------------------------------------------------------------
# 2) Spanish paragraph: split the single run describing the system into five
#    runs, inserting the new "por ultimo..." sentence and tweaking the tail.
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("A diferencia de estos sistemas, el proyecto que se presenta en este artículo está orientado exclusivamente a la gestión administrativa por parte del personal autorizado (administradores y bibliotecarios), excluyendo a los lectores de la manipulación directa del sistema. Además, incluye características distintivas como una interfaz intuitiva, un módulo de gestión de sanciones o multas por retraso o daños, y la gestión de cubículos de estudio disponibles en la biblioteca. Estas funcionalidades amplían el alcance del sistema, adaptándolo a las necesidades reales de los entornos bibliotecarios actuales.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sysPara = $rng2.Paragraphs(1)

$rPr = '<w:rPr><w:spacing w:val="0"/><w:lang w:val="es-EC"/></w:rPr>'
$run1 = '<w:r w:rsidRPr="00B9598A">' + $rPr + '<w:t>A diferencia de estos sistemas, el proyecto que se presenta en este artículo está orientado exclusivamente a la gestión administrativa por parte del personal autorizado (administradores y bibliotecarios), excluyendo a los lectores de la manipulación directa del sistema. Además, incluye características distintivas como una interfaz intuitiva, un módulo de gestión de sanciones o multas por retraso o daños, y la gestión de cubículos de estudio disponibles en la biblioteca</w:t></w:r>'
$run2 = '<w:r>' + $rPr + '<w:t>, por último, el proyecto realizado no cuenta con un sistema de facturación debido a que la biblioteca ya cuenta con uno propio y no quiere cambiarlo, siendo este proyecto por así decirlo un sistema administrativo para visualizar cada parte importante de la biblioteca sin recurrir a los libros que se tenía anteriormente</w:t></w:r>'
$run3 = '<w:r>' + $rPr + '<w:t>. Estas funcionalidades amplían el alcance del sistema, adaptándolo a las necesidades reales de</w:t></w:r>'
$run4 = '<w:r>' + $rPr + '<w:t xml:space="preserve"> la biblioteca</w:t></w:r>'
$run5 = '<w:r>' + $rPr + '<w:t>.</w:t></w:r>'
$sysRunsXml = $run1 + $run2 + $run3 + $run4 + $run5
Set-ParagraphRunsXml $sysPara $sysRunsXml

# ---------------------------------------------------------------------------
# 3) Heading: "2.2 Componentes del Sistema" -> "Componentes del Sistema"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("2.2 Componentes del Sistema", $true, $false, $false, $false, $false, $true, 1, $false, "Componentes del Sistema", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Reference list entry: collapse the COCOMO Model citation run group into
#    a single run (dropping the proofErr markers / extra nbsp-only run).
# ---------------------------------------------------------------------------
$rng4 = $d.Content
$rng4.Find.Execute([char]0x00AB + "Software Engineering | COCOMO Model," + [char]0x00BB, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
Write-Output ("cocomo find result text: [" + $rng4.Text + "]")
